$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = 111964632
$ws.Range("B2").Value2 = 77515
$ws.Range("D2").Value2 = "NT"
$ws.Range("E2").Value2 = 6425
$ws.Range("F2").Value2 = "Garnlav"
$ws.Range("G2").Value2 = "Alectoria sarmentosa"
$ws.Range("H2").Value2 = "(Ach.) Ach."
$ws.Range("Q2").Value2 = 734972.3834676194
$ws.Range("R2").Value2 = 7088252.533270728
$ws.Range("Z2").Value2 = "16:12"
$ws.Range("AB2").Value2 = "16:12"

# Row 3
$ws.Range("A3").Value2 = 111964847
$ws.Range("B3").Value2 = 89405
$ws.Range("D3").Value2 = "NT"
$ws.Range("E3").Value2 = 1202
$ws.Range("F3").Value2 = "Ullticka"
$ws.Range("G3").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("H3").Value2 = "(P.Karst.) Fiasson & Niemelä"

# Row 4
$ws.Range("A4").Value2 = 111964457
$ws.Range("B4").Value2 = 56398
$ws.Range("D4").Value2 = "NT"
$ws.Range("E4").Value2 = 100109
$ws.Range("F4").Value2 = "Tretåig hackspett"
$ws.Range("G4").Value2 = "Picoides tridactylus"
$ws.Range("H4").Value2 = "(Linnaeus, 1758)"
$ws.Range("Q4").Value2 = 734949.4564622594
$ws.Range("R4").Value2 = 7088268.525185317
$ws.Range("Z4").Value2 = "16:01"
$ws.Range("AB4").Value2 = "16:01"

# Row 5
$ws.Range("A5").Value2 = 111965883
$ws.Range("B5").Value2 = 55611
$ws.Range("E5").Value2 = 102612
$ws.Range("F5").Value2 = "Järpe"
$ws.Range("G5").Value2 = "Tetrastes bonasia"
$ws.Range("H5").Value2 = "(Linnaeus, 1758)"
$ws.Range("M5").Value2 = "lockläte, övriga läten"
$ws.Range("Q5").Value2 = 734846.6442297549
$ws.Range("R5").Value2 = 7088238.22626837
$ws.Range("Z5").Value2 = "17:05"
$ws.Range("AB5").Value2 = "17:05"

# Row 6
$ws.Range("A6").Value2 = 111964175
$ws.Range("B6").Value2 = 89423
$ws.Range("E6").Value2 = 5432
$ws.Range("F6").Value2 = "Granticka"
$ws.Range("G6").Value2 = "Porodaedalea chrysoloma"
$ws.Range("H6").Value2 = "(Fr.) Fiasson & Niemelä"
$ws.Range("M6").Value2 = ""
$ws.Range("Q6").Value2 = 734896.4627943118
$ws.Range("R6").Value2 = 7088342.483217424
$ws.Range("Z6").Value2 = "15:42"
$ws.Range("AB6").Value2 = "15:42"

# Row 7
$ws.Range("A7").Value2 = 111965370
$ws.Range("B7").Value2 = 81248
$ws.Range("E7").Value2 = 1312
$ws.Range("F7").Value2 = "Gammelgransskål"
$ws.Range("G7").Value2 = "Pseudographis pinicola"
$ws.Range("H7").Value2 = "(Nyl.) Rehm"
$ws.Range("Q7").Value2 = 734939.7547518623
$ws.Range("R7").Value2 = 7088232.371273324
$ws.Range("Z7").Value2 = "16:38"
$ws.Range("AB7").Value2 = "16:38"

# Row 8
$ws.Range("A8").Value2 = 111964050
$ws.Range("B8").Value2 = 90065
$ws.Range("D8").Value2 = "VU"
$ws.Range("E8").Value2 = 898
$ws.Range("F8").Value2 = "Blackticka"
$ws.Range("G8").Value2 = "Steccherinum collabens"
$ws.Range("H8").Value2 = "(Fr.) Vesterholt"
$ws.Range("Q8").Value2 = 734893.3330648565
$ws.Range("R8").Value2 = 7088354.646951701
$ws.Range("Z8").Value2 = "15:42"
$ws.Range("AB8").Value2 = "15:42"

# Row 9
$ws.Range("A9").Value2 = 111964863
$ws.Range("B9").Value2 = 89745
$ws.Range("D9").Value2 = "VU"
$ws.Range("E9").Value2 = 2062
$ws.Range("F9").Value2 = "Ulltickeporing"
$ws.Range("G9").Value2 = "Skeletocutis brevispora"
$ws.Range("H9").Value2 = "Niemelä"
$ws.Range("Q9").Value2 = 734972.3834676194
$ws.Range("R9").Value2 = 7088252.533270728
$ws.Range("Z9").Value2 = "16:12"
$ws.Range("AB9").Value2 = "16:12"

# Row 10
$ws.Range("A10").Value2 = 111965439
$ws.Range("B10").Value2 = 56398
$ws.Range("E10").Value2 = 100109
$ws.Range("F10").Value2 = "Tretåig hackspett"
$ws.Range("G10").Value2 = "Picoides tridactylus"
$ws.Range("H10").Value2 = "(Linnaeus, 1758)"
$ws.Range("Q10").Value2 = 734926.7697699566
$ws.Range("R10").Value2 = 7088234.05367971
$ws.Range("Z10").Value2 = "16:40"
$ws.Range("AB10").Value2 = "16:40"

# Row 11
$ws.Range("A11").Value2 = 111964622
$ws.Range("B11").Value2 = 89845
$ws.Range("D11").Value2 = "VU"
$ws.Range("E11").Value2 = 1209
$ws.Range("F11").Value2 = "Rynkskinn"
$ws.Range("G11").Value2 = "Phlebia centrifuga"
$ws.Range("H11").Value2 = "P.Karst."
$ws.Range("M11").Value2 = ""
$ws.Range("Q11").Value2 = 734972.3834676194
$ws.Range("R11").Value2 = 7088252.533270728
$ws.Range("Z11").Value2 = "16:12"
$ws.Range("AB11").Value2 = "16:12"
